$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column S: header timestamp, styled like the other header cells (copy R1's format).
$ws.Range("S1").Value2 = "2026-01-28 10:17:26"
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)

# Data rows 2-100: duplicate the latest price snapshot (column R) into the new column S.
for ($r = 2; $r -le 100; $r++) {
    $ws.Cells.Item($r, 19).Value2 = $ws.Cells.Item($r, 18).Value2
}

# Rows 101-204 have no price history yet; keep column S an empty text cell like D:R there.
for ($r = 101; $r -le 204; $r++) {
    $ws.Cells.Item($r, 19).Value2 = "'"
    $ws.Cells.Item($r, 19).Style = "Normal"
}
